# RD-Release history.xlsx — add the 7.0.3 release row, resize column K,
# and move the active selection, per the authored commit:
#   "The release FW for 7.0.3. ... Update the release note and the BT FW for V3."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New release-history row for firmware 7.0.3 (row 18 was blank, matching the
# pattern already used by rows 16 (7.0.1) and 17 (7.0.2)).
$ws.Range("A18").Value = "V3 EVT Firmware"
$ws.Range("B18").Value = "7.0.3"
$ws.Range("C18").Value = "3/31/2022"
$ws.Range("D18").Value = "Zound_Hendrix_M_Lite_V3_hwEVT_btswv7.0.3_20220331"
$ws.Range("E18").Value = "7.0.3"
$ws.Range("G18").Value = 0.6
$ws.Range("H18").Value = 3.1
$ws.Range("J18").Value = "N/A"
$ws.Range("K18").Value = "BT output gain correct. Tone play flow corrected."

# Widen the Comments column (K) so the new, longer note fits.
$ws.Columns.Item(11).ColumnWidth = 46.8

# Leave the cursor on the newly-added comment cell.
$ws.Range("K20").Select()
